$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (values scraped on 2024-02-11).
# D-column "Price" cells are plain text in the source sheet (some contain
# thousand-separator dots like "48.359.61" that are not valid numbers, and
# others like "40.50"/"2.00" would lose their trailing zero if Excel
# auto-coerced them to a Number). Force text storage explicitly, then reset
# the cell style back to Normal so no stray "Text" number-format style is
# left attached to the cell (keeps styles identical to the original file).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '48.359.61'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.66%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.525.95'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.90%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '322.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.84%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '109.39'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.37%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.533'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.33%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("E9").Value = '  +4.01%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '40.50'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.60%  '

$ws.Range("E11").Value = '  +13.69%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0825'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.31%  '

$ws.Range("E13").Value = '  +1.32%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.30'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.78%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.919.85'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.62%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.530.59'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.26%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '48.170.70'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.42%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.53'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +6.97%  '

$ws.Range("E20").Value = '  +0.76%  '

$ws.Range("E21").Value = '  +2.12%  '

$ws.Range("E22").Value = '  -0.91%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.64%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '264.73'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +8.28%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.58'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.21%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '26.12'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.09%  '

$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.30'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.81%  '

$ws.Range("B29").Value = 'Cosmos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.12'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.35%  '

$ws.Range("E30").Value = '  +8.24%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '35.73'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.66%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.69'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.57%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.81'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.38%  '

$ws.Range("E34").Value = '  +1.34%  '

$ws.Range("E35").Value = '  +0.00%  '

$ws.Range("E36").Value = '  +1.11%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.00'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.14%  '

$ws.Range("E38").Value = '  +2.24%  '

$ws.Range("E39").Value = '  +2.79%  '

$ws.Range("E40").Value = '  +0.85%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '120.91'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.30%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.17'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.52%  '

$ws.Range("E43").Value = '  -0.62%  '

$ws.Range("E44").Value = '  +2.68%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.018.36'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.11%  '

$ws.Range("E46").Value = '  +5.89%  '

$ws.Range("E47").Value = '  +8.51%  '

$ws.Range("E48").Value = '  +1.95%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.16'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.21%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '5.22'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.92%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '79.28'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.46%  '
